$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# C3 (Ns/m) and C4 (Ns/m) values changed from 100 to 0
$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 0

# torsional stiffness (kNm/rad) formula replaced with static value 12500
$ws.Range("B15").Value = 12500

# update selection to P9
$ws.Range("P9").Select()
